$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "68.665.32"
$ws.Cells.Item(2, 5).Value = "  -1.03%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.859.90"
$ws.Cells.Item(3, 5).Value = "  -2.14%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.09%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'522.46"
$ws.Cells.Item(5, 5).Value = "  +5.80%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'141.04"
$ws.Cells.Item(6, 5).Value = "  -4.74%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.609"
$ws.Cells.Item(7, 5).Value = "  -2.29%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 5).Value = "  +0.19%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.712"
$ws.Cells.Item(9, 5).Value = "  -3.07%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -5.68%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0000322"
$ws.Cells.Item(11, 5).Value = "  -8.14%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'41.67"
$ws.Cells.Item(12, 5).Value = "  -3.85%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'10.36"
$ws.Cells.Item(13, 5).Value = "  -0.79%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "4.491.33"
$ws.Cells.Item(14, 5).Value = "  -1.84%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'21.42"
$ws.Cells.Item(15, 5).Value = "  +7.82%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "3.875.14"
$ws.Cells.Item(16, 5).Value = "  -2.47%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "Uniswap"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(17, 4).Value = "'14.08"
$ws.Cells.Item(17, 5).Value = "  -2.06%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -2.16%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'1.19"
$ws.Cells.Item(19, 5).Value = "  +2.15%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "68.672.67"
$ws.Cells.Item(20, 5).Value = "  -1.10%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'416.61"
$ws.Cells.Item(21, 5).Value = "  -5.42%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +1.95%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'14.03"
$ws.Cells.Item(23, 5).Value = "  -3.29%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'86.88"
$ws.Cells.Item(24, 5).Value = "  -2.31%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +6.24%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'11.86"
$ws.Cells.Item(26, 5).Value = "  -2.18%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'10.45"
$ws.Cells.Item(27, 5).Value = "  -6.02%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'35.45"
$ws.Cells.Item(28, 5).Value = "  -4.62%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'13.47"
$ws.Cells.Item(29, 5).Value = "  +1.04%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'675.56"
$ws.Cells.Item(30, 5).Value = "  -4.57%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).Value = "'0.125"
$ws.Cells.Item(31, 5).Value = "  -5.13%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "NEARProtocol"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(32, 4).Value = "'6.90"
$ws.Cells.Item(32, 5).Value = "  +13.36%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'2.84"
$ws.Cells.Item(33, 5).Value = "  -1.79%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'66.76"
$ws.Cells.Item(34, 5).Value = "  +8.46%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.445"
$ws.Cells.Item(35, 5).Value = "  -5.15%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "0.0₃0852"
$ws.Cells.Item(36, 5).Value = "  -7.31%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'39.45"
$ws.Cells.Item(37, 5).Value = "  -3.52%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'3.51"
$ws.Cells.Item(38, 5).Value = "  +15.17%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -1.69%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.998"
$ws.Cells.Item(40, 5).Value = "  +0.09%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'1.00"
$ws.Cells.Item(41, 5).Value = "  -0.11%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'0.0474"
$ws.Cells.Item(42, 5).Value = "  -3.32%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Fetch.AI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(43, 4).Value = "'2.86"
$ws.Cells.Item(43, 5).Value = "  -2.92%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "WEMIXToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(44, 4).Value = "'3.16"
$ws.Cells.Item(44, 5).Value = "  +4.83%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'3.41"
$ws.Cells.Item(45, 5).Value = "  +1.90%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "FLOKI"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(46, 4).Value = "'0.000288"
$ws.Cells.Item(46, 5).Value = "  +18.22%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Stellar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(47, 4).Value = "'0.141"
$ws.Cells.Item(47, 5).Value = "  -1.94%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Stacks"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(48, 4).Value = "'2.98"
$ws.Cells.Item(48, 5).Value = "  -1.77%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'3.28"
$ws.Cells.Item(49, 5).Value = "  -3.19%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'8.77"
$ws.Cells.Item(50, 5).Value = "  +3.88%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'142.79"
$ws.Cells.Item(51, 5).Value = "  -0.86%  "
